{"js": "// Replace \"EC2 instances\" with \"Amazon EC2 instances\" everywhere in the\n// document body (the four identical \"System Maintenance control ...\"\n// paragraphs that mention EC2 instances).\nconst body = context.document.body;\nconst results = body.search(\"EC2 instances\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Amazon EC2 instances\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace \"EC2 instances\" with \"Amazon EC2 instances\" everywhere in the\n# document (the four identical \"System Maintenance control ...\" paragraphs\n# that mention EC2 instances).\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"EC2 instances\"\n$find.Replacement.Text = \"Amazon EC2 instances\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute(\n  [ref]$find.Text,\n  [ref]$false,\n  [ref]$false,\n  [ref]$false,\n  [ref]$false,\n  [ref]$false,\n  [ref]$true,\n  [ref]$wdFindContinue,\n  [ref]$false,\n  [ref]$find.Replacement.Text,\n  [ref]$wdReplaceAll\n)\n"}
